$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 329, pushing the existing rows 329-333
# down to 331-335 (formatting/formulas are carried down automatically).
$ws.Rows("329:330").Insert()

# --- New row 329: "Primera" quality entry for 2022-02-08 (serial 44595) ---
$ws.Cells.Item(329, 1).Value2 = 4
$ws.Cells.Item(329, 2).Value2 = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(329, 3).Value2 = "Los Lagos"
$ws.Cells.Item(329, 4).Value2 = 44595
$ws.Cells.Item(329, 5).Value2 = 10
$ws.Cells.Item(329, 6).Value2 = "Fruta"
$ws.Cells.Item(329, 7).Value2 = 100106
$ws.Cells.Item(329, 8).Value2 = "Oleaginosos"
$ws.Cells.Item(329, 9).Value2 = 100106002
$ws.Cells.Item(329, 10).Value2 = "Palta"
$ws.Cells.Item(329, 11).Value2 = "Hass"
$ws.Cells.Item(329, 12).Value2 = "Primera"
$ws.Cells.Item(329, 13).Value2 = 200
$ws.Cells.Item(329, 14).Value2 = 4000
$ws.Cells.Item(329, 15).Value2 = 4100
$ws.Cells.Item(329, 16).Value2 = 4050
$ws.Cells.Item(329, 17).Value2 = "`$/kilo (en caja de 17 kilos)"
$ws.Cells.Item(329, 18).Value2 = "Provincia de Quillota"
$ws.Cells.Item(329, 19).Value2 = 4050
$ws.Cells.Item(329, 20).Value2 = 1

# --- New row 330: "Segunda" quality entry for 2022-02-08 (serial 44595) ---
$ws.Cells.Item(330, 1).Value2 = 4
$ws.Cells.Item(330, 2).Value2 = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(330, 3).Value2 = "Los Lagos"
$ws.Cells.Item(330, 4).Value2 = 44595
$ws.Cells.Item(330, 5).Value2 = 10
$ws.Cells.Item(330, 6).Value2 = "Fruta"
$ws.Cells.Item(330, 7).Value2 = 100106
$ws.Cells.Item(330, 8).Value2 = "Oleaginosos"
$ws.Cells.Item(330, 9).Value2 = 100106002
$ws.Cells.Item(330, 10).Value2 = "Palta"
$ws.Cells.Item(330, 11).Value2 = "Hass"
$ws.Cells.Item(330, 12).Value2 = "Segunda"
$ws.Cells.Item(330, 13).Value2 = 100
$ws.Cells.Item(330, 14).Value2 = 3500
$ws.Cells.Item(330, 15).Value2 = 3500
$ws.Cells.Item(330, 16).Value2 = 3500
$ws.Cells.Item(330, 17).Value2 = "`$/kilo (en caja de 17 kilos)"
$ws.Cells.Item(330, 18).Value2 = "Provincia de Quillota"
$ws.Cells.Item(330, 19).Value2 = 3500
$ws.Cells.Item(330, 20).Value2 = 1
